$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Range B2:F25
$valuesBF = New-Object 'object[,]' 24,5
$valuesBF[0,0] = 1.02
$valuesBF[0,1] = 1.040970137894234
$valuesBF[0,2] = 1.035769314310167
$valuesBF[0,3] = 1.039240551556146
$valuesBF[0,4] = 1.048477071803056
$valuesBF[1,0] = 1.02
$valuesBF[1,1] = 1.042428100992114
$valuesBF[1,2] = 1.036401314810634
$valuesBF[1,3] = 1.040497502991863
$valuesBF[1,4] = 1.050049162337545
$valuesBF[2,0] = 1.02
$valuesBF[2,1] = 1.043370028331797
$valuesBF[2,2] = 1.036809845544632
$valuesBF[2,3] = 1.041309760673375
$valuesBF[2,4] = 1.051065266594352
$valuesBF[3,0] = 1.02
$valuesBF[3,1] = 1.043765669365564
$valuesBF[3,2] = 1.036981492982147
$valuesBF[3,3] = 1.04165098201017
$valuesBF[3,4] = 1.051492170855837
$valuesBF[4,0] = 1.02
$valuesBF[4,1] = 1.043832079202001
$valuesBF[4,2] = 1.037010307562834
$valuesBF[4,3] = 1.041708259962601
$valuesBF[4,4] = 1.051563834588448
$valuesBF[5,0] = 1.02
$valuesBF[5,1] = 1.043375316254171
$valuesBF[5,2] = 1.036812139493839
$valuesBF[5,3] = 1.041314321072387
$valuesBF[5,4] = 1.051070971945442
$valuesBF[6,0] = 1.02
$valuesBF[6,1] = 1.041463172566681
$valuesBF[6,2] = 1.035982988353406
$valuesBF[6,3] = 1.039665569754117
$valuesBF[6,4] = 1.049008608608218
$valuesBF[7,0] = 1.02
$valuesBF[7,1] = 1.038082129536699
$valuesBF[7,2] = 1.034518701845919
$valuesBF[7,3] = 1.036751795036302
$valuesBF[7,4] = 1.045365354179897
$valuesBF[8,0] = 1.02
$valuesBF[8,1] = 1.035819843178793
$valuesBF[8,2] = 1.033540301476733
$valuesBF[8,3] = 1.034803250961291
$valuesBF[8,4] = 1.042929904579893
$valuesBF[9,0] = 1.02
$valuesBF[9,1] = 1.034838187197287
$valuesBF[9,2] = 1.033116107893378
$valuesBF[9,3] = 1.033958004551284
$valuesBF[9,4] = 1.041873649698652
$valuesBF[10,0] = 1.02
$valuesBF[10,1] = 1.034473236400098
$valuesBF[10,2] = 1.032958461338656
$valuesBF[10,3] = 1.033643808339312
$valuesBF[10,4] = 1.041481046303338
$valuesBF[11,0] = 1.02
$valuesBF[11,1] = 1.034551534211914
$valuesBF[11,2] = 1.032992280795736
$valuesBF[11,3] = 1.033711215191702
$valuesBF[11,4] = 1.04156527313314
$valuesBF[12,0] = 1.02
$valuesBF[12,1] = 1.034808026831369
$valuesBF[12,2] = 1.033103078462563
$valuesBF[12,3] = 1.033932037807553
$valuesBF[12,4] = 1.041841202406202
$valuesBF[13,0] = 1.02
$valuesBF[13,1] = 1.034966017721287
$valuesBF[13,2] = 1.033171333640827
$valuesBF[13,3] = 1.03406806269543
$valuesBF[13,4] = 1.042011176361664
$valuesBF[14,0] = 1.02
$valuesBF[14,1] = 1.035884948059973
$valuesBF[14,2] = 1.033568442344209
$valuesBF[14,3] = 1.034859314721759
$valuesBF[14,4] = 1.042999968337844
$valuesBF[15,0] = 1.02
$valuesBF[15,1] = 1.036460808037008
$valuesBF[15,2] = 1.033817392878802
$valuesBF[15,3] = 1.035355236235519
$valuesBF[15,4] = 1.043619752486899
$valuesBF[16,0] = 1.02
$valuesBF[16,1] = 1.036796498197095
$valuesBF[16,2] = 1.033962549508954
$valuesBF[16,3] = 1.035644353237879
$valuesBF[16,4] = 1.043981099766248
$valuesBF[17,0] = 1.02
$valuesBF[17,1] = 1.036910926343296
$valuesBF[17,2] = 1.034012035379609
$valuesBF[17,3] = 1.035742910162168
$valuesBF[17,4] = 1.044104282603602
$valuesBF[18,0] = 1.02
$valuesBF[18,1] = 1.036399044382321
$valuesBF[18,2] = 1.033790688229455
$valuesBF[18,3] = 1.035302043661149
$valuesBF[18,4] = 1.043553272372773
$valuesBF[19,0] = 1.02
$valuesBF[19,1] = 1.034732505054849
$valuesBF[19,2] = 1.033070453588325
$valuesBF[19,3] = 1.033867017565574
$valuesBF[19,4] = 1.041759955457978
$valuesBF[20,0] = 1.02
$valuesBF[20,1] = 1.033682829421452
$valuesBF[20,2] = 1.032617137517595
$valuesBF[20,3] = 1.032963402189435
$valuesBF[20,4] = 1.040630897775974
$valuesBF[21,0] = 1.02
$valuesBF[21,1] = 1.034239461272559
$valuesBF[21,2] = 1.032857494425407
$valuesBF[21,3] = 1.03344255669232
$valuesBF[21,4] = 1.041229580667423
$valuesBF[22,0] = 1.02
$valuesBF[22,1] = 1.036426953335096
$valuesBF[22,2] = 1.033802755071998
$valuesBF[22,3] = 1.035326079546102
$valuesBF[22,4] = 1.043583312376557
$valuesBF[23,0] = 1.02
$valuesBF[23,1] = 1.038957631941373
$valuesBF[23,2] = 1.034897640897746
$valuesBF[23,3] = 1.037506113202995
$valuesBF[23,4] = 1.046308352255047
$ws.Range("B2:F25").Value = $valuesBF

# Range I2:N25
$valuesIN = New-Object 'object[,]' 24,6
$valuesIN[0,0] = 1.038624342265089
$valuesIN[0,1] = 1.046053871020944
$valuesIN[0,2] = 1.038564970202844
$valuesIN[0,3] = 1.042026301067575
$valuesIN[0,4] = 1.051236808243076
$valuesIN[0,5] = 1.04753938750185
$valuesIN[1,0] = 1.038969925624626
$valuesIN[1,1] = 1.047155136229404
$valuesIN[1,2] = 1.039006905657447
$valuesIN[1,3] = 1.043092252588891
$valuesIN[1,4] = 1.052618986547395
$valuesIN[1,5] = 1.048642216633223
$valuesIN[2,0] = 1.039191802847452
$valuesIN[2,1] = 1.047865862237047
$valuesIN[2,2] = 1.039291763949868
$valuesIN[2,3] = 1.043780371278091
$valuesIN[2,4] = 1.053511720426428
$valuesIN[2,5] = 1.049353951953318
$valuesIN[3,0] = 1.039284665481417
$valuesIN[3,1] = 1.048164209975367
$valuesIN[3,2] = 1.039411255265331
$valuesIN[3,3] = 1.044069272671989
$valuesIN[3,4] = 1.053886643002361
$valuesIN[3,5] = 1.049652723379647
$valuesIN[4,0] = 1.03930023326328
$valuesIN[4,1] = 1.048214278138555
$valuesIN[4,2] = 1.039431302984838
$valuesIN[4,3] = 1.044117758128529
$valuesIN[4,4] = 1.053949571910175
$valuesIN[4,5] = 1.049702862645371
$valuesIN[5,0] = 1.039193045309435
$valuesIN[5,1] = 1.047869850502283
$valuesIN[5,2] = 1.039293361631217
$valuesIN[5,3] = 1.043784233094302
$valuesIN[5,4] = 1.053516731655549
$valuesIN[5,5] = 1.049357945882348
$valuesIN[6,0] = 1.038741494680902
$valuesIN[6,1] = 1.046426438368638
$valuesIN[6,2] = 1.038714553493386
$valuesIN[6,3] = 1.042386883432363
$valuesIN[6,4] = 1.051704263152423
$valuesIN[6,5] = 1.047912483937912
$valuesIN[7,0] = 1.037932418764973
$valuesIN[7,1] = 1.043868430166539
$valuesIN[7,2] = 1.037686119392562
$valuesIN[7,3] = 1.039911922310734
$valuesIN[7,4] = 1.048497680612506
$valuesIN[7,5] = 1.045350843070758
$valuesIN[8,0] = 1.037383931910538
$valuesIN[8,1] = 1.042152985318462
$valuesIN[8,2] = 1.036994711344854
$valuesIN[8,3] = 1.038253130835219
$valuesIN[8,4] = 1.0463509247771
$valuesIN[8,5] = 1.043632962094232
$valuesIN[9,0] = 1.037144248607788
$valuesIN[9,1] = 1.041407706344718
$valuesIN[9,2] = 1.036693935545322
$valuesIN[9,3] = 1.037532694667692
$valuesIN[9,4] = 1.04541911436813
$valuesIN[9,5] = 1.042886624738861
$valuesIN[10,0] = 1.037054889228102
$valuesIN[10,1] = 1.041130497367826
$valuesIN[10,2] = 1.036582003478058
$valuesIN[10,3] = 1.037264760705029
$valuesIN[10,4] = 1.045072651729154
$valuesIN[10,5] = 1.042609022093426
$valuesIN[11,0] = 1.037074072100018
$valuesIN[11,1] = 1.041189976899487
$valuesIN[11,2] = 1.036606022837931
$valuesIN[11,3] = 1.037322248566829
$valuesIN[11,4] = 1.04514698497162
$valuesIN[11,5] = 1.042668586092844
$valuesIN[12,0] = 1.037136868884283
$valuesIN[12,1] = 1.041384799937706
$valuesIN[12,2] = 1.03668468750937
$valuesIN[12,3] = 1.03751055396872
$valuesIN[12,4] = 1.045390482782079
$valuesIN[12,5] = 1.042863685802124
$valuesIN[13,0] = 1.037175516221356
$valuesIN[13,1] = 1.041504786398242
$valuesIN[13,2] = 1.036733127466097
$valuesIN[13,3] = 1.037626530979405
$valuesIN[13,4] = 1.045540463612046
$valuesIN[13,5] = 1.042983842657196
$valuesIN[14,0] = 1.037399792692373
$valuesIN[14,1] = 1.042202394217083
$valuesIN[14,2] = 1.037014643399978
$valuesIN[14,3] = 1.03830089760353
$valuesIN[14,4] = 1.046412717808826
$valuesIN[14,5] = 1.043682441159157
$valuesIN[15,0] = 1.037539888984319
$valuesIN[15,1] = 1.042639316562991
$valuesIN[15,2] = 1.037190857312368
$valuesIN[15,3] = 1.038723325108238
$valuesIN[15,4] = 1.046959251136034
$valuesIN[15,5] = 1.044119983984912
$valuesIN[16,0] = 1.03762139411403
$valuesIN[16,1] = 1.04289392713603
$valuesIN[16,2] = 1.037293505743065
$valuesIN[16,3] = 1.03896951110993
$valuesIN[16,4] = 1.047277818177252
$valuesIN[16,5] = 1.044374956134168
$valuesIN[17,0] = 1.037649149586542
$valuesIN[17,1] = 1.042980702448598
$valuesIN[17,2] = 1.037328483493056
$valuesIN[17,3] = 1.039053418939331
$valuesIN[17,4] = 1.04738640482473
$valuesIN[17,5] = 1.044461854677634
$valuesIN[18,0] = 1.037524879783413
$valuesIN[18,1] = 1.042592463682263
$valuesIN[18,2] = 1.037171965115115
$valuesIN[18,3] = 1.038678024277714
$valuesIN[18,4] = 1.046900635748709
$valuesIN[18,5] = 1.044073064567719
$valuesIN[19,0] = 1.037118385931959
$valuesIN[19,1] = 1.041327439929232
$valuesIN[19,2] = 1.036661528548296
$valuesIN[19,3] = 1.037455111932336
$valuesIN[19,4] = 1.045318788387564
$valuesIN[19,5] = 1.042806244335859
$valuesIN[20,0] = 1.036860894888534
$valuesIN[20,1] = 1.04052987215331
$valuesIN[20,2] = 1.036339377933636
$valuesIN[20,3] = 1.036684295118138
$valuesIN[20,4] = 1.044322207598358
$valuesIN[20,5] = 1.042007543922212
$valuesIN[21,0] = 1.036997577717449
$valuesIN[21,1] = 1.040952888396599
$valuesIN[21,2] = 1.036510272104714
$valuesIN[21,3] = 1.037093104049316
$valuesIN[21,4] = 1.044850707391204
$valuesIN[21,5] = 1.042431160897087
$valuesIN[22,0] = 1.037531662446304
$valuesIN[22,1] = 1.042613635220517
$valuesIN[22,2] = 1.037180502100595
$valuesIN[22,3] = 1.038698494418521
$valuesIN[22,4] = 1.04692712218867
$valuesIN[22,5] = 1.044094266171986
$valuesIN[23,0] = 1.03814318112776
$valuesIN[23,1] = 1.044531493084362
$valuesIN[23,2] = 1.03795300874248
$valuesIN[23,3] = 1.040553290293439
$valuesIN[23,4] = 1.049328218896717
$valuesIN[23,5] = 1.046014847613979
$ws.Range("I2:N25").Value = $valuesIN
